# refs #684 Risiko 4 bereinigt
# Add the new changelog entry (version 1.7) to the "Aenderungsgeschichte" table
# and clear out Risiko 4 ("Kinect: Menschliche Drehung") in the "Risiken" sheet.

$wb = $excel.ActiveWorkbook
$wsChangelog = $wb.Worksheets.Item(1)
$wsRisiken   = $wb.Worksheets.Item(2)

# --- 1. Grow the changelog table (Table2 on sheet 1) by one row ------------
$changelogTable = $wsChangelog.ListObjects.Item(1)
$changelogTable.ListRows.Add() | Out-Null

$wsChangelog.Cells.Item(11, 1).Value = 41001
$wsChangelog.Cells.Item(11, 2).Value = "1.7"
$wsChangelog.Cells.Item(11, 3).Value = "Risiko 4 ""Kinect: Menschliche Drehung"" ist bereinigt. Zum jetztigen Projektzeitpunkt ist der Teaser, fuer welchen die Erkennung der Drehung relevant gewesen waere, tief priorisiert. Er wird nicht mehr innerhalb des Zeitraums der BA implementiert werden koennen. "
$wsChangelog.Cells.Item(11, 4).Value = "DT"

$a11 = $wsChangelog.Range("A11")
$a11.HorizontalAlignment = -4131
$a11.WrapText = $true
$a11.NumberFormat = "mm-dd-yy"

$b11 = $wsChangelog.Range("B11")
$b11.HorizontalAlignment = -4131
$b11.WrapText = $true
$b11.NumberFormat = "@"

$c11 = $wsChangelog.Range("C11")
$c11.HorizontalAlignment = -4131
$c11.WrapText = $true
$c11.VerticalAlignment = -4160

$d11 = $wsChangelog.Range("D11")
$d11.WrapText = $true

$wsChangelog.Rows.Item(11).RowHeight = 60

# --- 2. Risiko 4 "Kinect: Menschliche Drehung" bereinigt -------------------
$wsRisiken.Range("D7").Value = 0
$wsRisiken.Range("J7").Value = 6

# --- 3. Selections + active sheet ------------------------------------------
$wsChangelog.Range("D11").Select()
$wsRisiken.Range("J8").Select()
$wsRisiken.Activate()
